$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 87, shifting existing rows 87:144 down to 88:145.
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new weekly price entry.
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 44529
$ws.Range("D87").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112039
$ws.Range("G87").Value = "Ciboulette"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 80
$ws.Range("K87").Value = 2500
$ws.Range("L87").Value = 2500
$ws.Range("M87").Value = 2500
$ws.Range("N87").Value = "$/docena de atados"
$ws.Range("O87").Value = "Región Metropolitana"
$ws.Range("P87").Value = 833
$ws.Range("Q87").Value = 3
$ws.Range("R87").Value = "Hortaliza"
